$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (kept as text to match source formatting)
$updates = @{
    "D2" = "315.16"
    "E2" = "2.86%"
    "D3" = "35.26"
    "E3" = "-2.44%"
    "D4" = "5.134"
    "E4" = "0.42%"
    "D5" = "0.08127"
    "E5" = "3.22%"
    "D6" = "2.134"
    "E6" = "-0.06%"
    "D7" = "7.988"
    "E7" = "0.25%"
    "D8" = "4.149"
    "E8" = "0.61%"
    "D9" = "0.9297"
    "E9" = "1.19%"
    "D10" = "0.1007"
    "E10" = "3.72%"
    "E11" = "0.79%"
    "D12" = "0.09138"
    "E12" = "5.17%"
    "D13" = "0.03601"
    "E13" = "0.81%"
    "E14" = "-0.26%"
    "D15" = "0.001439"
    "E15" = "0.14%"
    "D16" = "0.005714"
    "E16" = "-0.39%"
    "D17" = "3.466"
    "E17" = "-0.10%"
    "E18" = "-1.05%"
    "E19" = "0.46%"
    "E20" = "0.07%"
    "D21" = "5.094"
    "E21" = "-1.28%"
    "E22" = "9.77%"
    "D23" = "0.04556"
    "E23" = "-0.22%"
    "D24" = "0.001245"
    "E24" = "0.54%"
    "D25" = "0.004707"
    "E25" = "-6.57%"
    "D26" = "0.0001251"
    "E26" = "-21.95%"
    "D27" = "0.0004503"
    "E27" = "-5.30%"
    "D39" = "0.01956"
    "E39" = "5.84%"
    "D40" = "0.04845"
    "E40" = "1.78%"
    "D41" = "0.007717"
    "E41" = "1.64%"
    "E42" = "-0.62%"
    "D43" = "0.007844"
    "E43" = "1.16%"
    "D44" = "0.002127"
    "E44" = "-3.18%"
    "D45" = "0.01180"
    "E45" = "6.83%"
    "D46" = "0.00006620"
    "E46" = "4.56%"
    "E47" = "-0.05%"
    "D48" = "37.54"
    "E48" = "-21.11%"
    "D49" = "0.001701"
    "E49" = "-15.04%"
    "E50" = "-0.05%"
    "E51" = "-0.05%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
